$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row (row 1): "<name>_old" -> "<name>_FV2404"
#    and "<name>_new" -> "<name>_FV2410". Column K ("diff") is untouched.
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")
$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$rightCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($leftCols[$i]  + "1").Value = $baseNames[$i] + "_FV2404"
    $ws.Range($rightCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into an Excel Table ("Table1") without picking up
#    a table style override for the (already custom-styled) header row.
#    We temporarily stash the header's existing formatting, clear it so the
#    new ListObject does not snapshot it into a headerRowDxfId, recreate the
#    table, then restore the original header formatting.
# ---------------------------------------------------------------------------
$header  = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")

$header.Copy()
$scratch.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$header.ClearFormats()

$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U94"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""

$scratch.Copy()
$header.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$scratch.EntireRow.Delete()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "done"
